$wb = $excel.ActiveWorkbook
$srcSheet = $wb.Worksheets.Item("Bus_Makhulu_f")
$srcSheet.Copy($wb.Worksheets.Item(1))
$newSheet = $wb.Worksheets.Item(1)
Write-Host "New sheet at front:" $newSheet.Name
